# Fix: "erro funcional nao salvava os dados da carga no espelho do sal"
# Adds the missing rows that should have been appended to "Programacao" and
# "Planilha" when a load was registered, and fixes the values on the
# "Descarga do Sal" mirror/spelho sheet that were not being saved correctly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Programacao": two new rows (5 and 6) with the load entry data.
# ---------------------------------------------------------------------
$wsProg = $wb.Worksheets.Item("Programacao")

$wsProg.Range("A5").NumberFormat = "@"
$wsProg.Range("A5").Value = "03/04/2025"
$wsProg.Range("B5").NumberFormat = "@"
$wsProg.Range("B5").Value = "05:27"
$wsProg.Range("C5").Value = "valdik antonio sa silva"
$wsProg.Range("D5").Value = "sda"
$wsProg.Range("E5").Value = "dsa/rte"
$wsProg.Range("F5").Value = "Selecione uma opção"
$wsProg.Range("G5").Value = 651
$wsProg.Range("H5").Value = "dsa"
$wsProg.Range("I5").Value = "das"
$wsProg.Range("J5").Value = "Selecione uma opção Selecione uma opção"
$wsProg.Range("K5").Value = "das"

$wsProg.Range("A6").Value = "das"
$wsProg.Range("B6").Value = "das"
$wsProg.Range("C6").Value = "das"
$wsProg.Range("D6").Value = "das"
$wsProg.Range("E6").Value = "das/dsa"
$wsProg.Range("F6").Value = "Selecione uma opção"
$wsProg.Range("G6").Value = 379
$wsProg.Range("H6").Value = "das"
$wsProg.Range("I6").Value = "das"
$wsProg.Range("J6").Value = "Selecione uma opção Selecione uma opção"
$wsProg.Range("K6").Value = "das"

# ---------------------------------------------------------------------
# Sheet "Planilha": four new rows (8-11) mirroring the load entries.
# ---------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("Planilha")

$wsPlan.Range("A8").Value = "ENTRADA"
$wsPlan.Range("B8").NumberFormat = "@"
$wsPlan.Range("B8").Value = "03/04/2025"
$wsPlan.Range("C8").Value = "dsa"
$wsPlan.Range("D8").Value = "das"
$wsPlan.Range("E8").Value = "das"
$wsPlan.Range("F8").Value = "Selecione uma opção"
$wsPlan.Range("G8").Value = "Selecione uma opção"
$wsPlan.Range("H8").Value = "Selecione uma opção"
$wsPlan.Range("I8").Value = "dsa"
$wsPlan.Range("J8").Value = "das"
$wsPlan.Range("K8").Value = 84965
$wsPlan.Range("L8").Value = "das"
$wsPlan.Range("M8").Value = "das"
$wsPlan.Range("N8").Value = 651

$wsPlan.Range("A9").Value = "ENTRADA"
$wsPlan.Range("B9").NumberFormat = "@"
$wsPlan.Range("B9").Value = "03/04/2025"
$wsPlan.Range("C9").Value = "dsa"
$wsPlan.Range("D9").Value = "das"
$wsPlan.Range("E9").Value = "das"
$wsPlan.Range("F9").Value = "Selecione uma opção"
$wsPlan.Range("G9").Value = "Selecione uma opção"
$wsPlan.Range("H9").Value = "Selecione uma opção"
$wsPlan.Range("I9").Value = "rte"
$wsPlan.Range("J9").NumberFormat = "@"
$wsPlan.Range("J9").Value = "98"
$wsPlan.Range("K9").Value = 8956
$wsPlan.Range("L9").Value = "fe"
$wsPlan.Range("M9").Value = "das"
$wsPlan.Range("N9").Value = 98465

$wsPlan.Range("A10").Value = "ENTRADA"
$wsPlan.Range("B10").Value = "das"
$wsPlan.Range("C10").Value = "das"
$wsPlan.Range("D10").Value = "das"
$wsPlan.Range("E10").Value = "das"
$wsPlan.Range("F10").Value = "Selecione uma opção"
$wsPlan.Range("G10").Value = "Selecione uma opção"
$wsPlan.Range("H10").Value = "Selecione uma opção"
$wsPlan.Range("I10").Value = "das"
$wsPlan.Range("J10").Value = "das"
$wsPlan.Range("K10").Value = 249
$wsPlan.Range("L10").Value = "das"
$wsPlan.Range("M10").Value = "das"
$wsPlan.Range("N10").Value = 379

$wsPlan.Range("A11").Value = "ENTRADA"
$wsPlan.Range("B11").Value = "das"
$wsPlan.Range("C11").Value = "das"
$wsPlan.Range("D11").Value = "das"
$wsPlan.Range("E11").Value = "das"
$wsPlan.Range("F11").Value = "Selecione uma opção"
$wsPlan.Range("G11").Value = "Selecione uma opção"
$wsPlan.Range("H11").Value = "Selecione uma opção"
$wsPlan.Range("I11").Value = "dsa"
$wsPlan.Range("J11").Value = "das"
$wsPlan.Range("K11").Value = 2767
$wsPlan.Range("L11").Value = "das"
$wsPlan.Range("M11").Value = "das"
$wsPlan.Range("N11").Value = 5786

# ---------------------------------------------------------------------
# Sheet "Descarga do Sal": the mirror/spelho that was not saving the
# load data correctly - correct all the placeholder/stale values.
# ---------------------------------------------------------------------
$wsSal = $wb.Worksheets.Item("Descarga do Sal")

$wsSal.Range("D8").Value = "das"
$wsSal.Range("K8").Value = "das"
$wsSal.Range("D10").Value = "das"
$wsSal.Range("D12").Value = "das"
$wsSal.Range("D14").Value = "das"
$wsSal.Range("K14").Value = "das"
$wsSal.Range("D16").Value = "das"
$wsSal.Range("D18").Value = "Selecione uma opção"
$wsSal.Range("D20").Value = "das"
$wsSal.Range("K20").Value = "das"
$wsSal.Range("P20").Value = 379
$wsSal.Range("D22").Value = "dsa"
$wsSal.Range("K22").Value = "das"
$wsSal.Range("P22").Value = 5786
$wsSal.Range("D26").Value = "Selecione uma opção"
$wsSal.Range("L26").Value = "das"
$wsSal.Range("D28").Value = "das"
$wsSal.Range("H28").Value = "das"
$wsSal.Range("K28").Value = 379
$wsSal.Range("O28").Value = 249
$wsSal.Range("D30").Value = "das"
$wsSal.Range("H30").Value = "dsa"
$wsSal.Range("K30").Value = 5786
$wsSal.Range("O30").Value = 2767

Write-Host "Applied load-mirror fix."
